# Fixed import of CDP to include cost of investment
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (old E:G -> F:H), pushing "Payment Date *",
# "Distribution from Sale" and "Folio No" one column to the right.
$ws.Columns.Item(5).Insert()

# New column header
$ws.Range("E1").Value = "Cost Of Investment *"

# New column values (Cost Of Investment, per row)
$ws.Range("E2").Value = 100000
$ws.Range("E3").Value = 200000
$ws.Range("E4").Value = 90000
$ws.Range("E5").Value = 100000
$ws.Range("E6").Value = 200000
$ws.Range("E7").Value = 200000

# Match number formatting of the neighboring "Amount *" column
$ws.Range("E2:E7").NumberFormat = $ws.Range("D2").NumberFormat()

# Set the new column's width to match "Amount *" column's width
$ws.Columns.Item(5).ColumnWidth = 18

# Leave selection where the editor ended up after entering the data
$ws.Range("E8").Select() | Out-Null
